$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text (no number auto-conversion), preserving default (no explicit) style
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextCell $ws.Range("D2") "67.646.97"
Set-TextCell $ws.Range("E2") "  -0.56%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.481.95"
Set-TextCell $ws.Range("E3") "  -1.33%  "

# Row 4
Set-TextCell $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextCell $ws.Range("D5") "592.22"
Set-TextCell $ws.Range("E5") "  -1.92%  "

# Row 6
Set-TextCell $ws.Range("D6") "179.58"
Set-TextCell $ws.Range("E6") "  -1.34%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.611"
Set-TextCell $ws.Range("E7") "  +1.97%  "

# Row 8
Set-TextCell $ws.Range("E8") "  +0.01%  "

# Row 9
Set-TextCell $ws.Range("D9") "3.480.80"
Set-TextCell $ws.Range("E9") "  -1.37%  "

# Row 10
Set-TextCell $ws.Range("E10") "  -1.76%  "

# Row 11
Set-TextCell $ws.Range("E11") "  -2.74%  "

# Row 12
Set-TextCell $ws.Range("E12") "  -3.09%  "

# Row 13
Set-TextCell $ws.Range("D13") "4.088.84"
Set-TextCell $ws.Range("E13") "  -1.24%  "

# Row 14
Set-TextCell $ws.Range("D14") "32.35"
Set-TextCell $ws.Range("E14") "  +0.22%  "

# Row 15
Set-TextCell $ws.Range("E15") "  -2.60%  "

# Row 16
Set-TextCell $ws.Range("D16") "67.630.18"
Set-TextCell $ws.Range("E16") "  -0.52%  "

# Row 17
Set-TextCell $ws.Range("E17") "  -2.44%  "

# Row 18
Set-TextCell $ws.Range("D18") "3.483.55"
Set-TextCell $ws.Range("E18") "  -1.41%  "

# Row 19
Set-TextCell $ws.Range("E19") "  -3.85%  "

# Row 20
Set-TextCell $ws.Range("D20") "14.06"
Set-TextCell $ws.Range("E20") "  -3.54%  "

# Row 21
Set-TextCell $ws.Range("D21") "388.38"
Set-TextCell $ws.Range("E21") "  -3.54%  "

# Row 22
Set-TextCell $ws.Range("D22") "7.91"
Set-TextCell $ws.Range("E22") "  -1.66%  "

# Row 23
Set-TextCell $ws.Range("E23") "  +1.54%  "

# Row 24
Set-TextCell $ws.Range("E24") "  -0.04%  "

# Row 25
Set-TextCell $ws.Range("D25") "72.36"
Set-TextCell $ws.Range("E25") "  -2.30%  "

# Row 26
Set-TextCell $ws.Range("E26") "  -1.84%  "

# Row 27
Set-TextCell $ws.Range("D27") "0.0000122"
Set-TextCell $ws.Range("E27") "  -1.15%  "

# Row 28
Set-TextCell $ws.Range("D28") "10.11"

# Row 29
Set-TextCell $ws.Range("E29") "  -1.57%  "

# Row 30
Set-TextCell $ws.Range("E30") "  +0.24%  "

# Row 31
Set-TextCell $ws.Range("D31") "6.06"
Set-TextCell $ws.Range("E31") "  -4.40%  "

# Row 32
Set-TextCell $ws.Range("D32") "24.66"
Set-TextCell $ws.Range("E32") "  +2.42%  "

# Row 33
Set-TextCell $ws.Range("D33") "2.04"
Set-TextCell $ws.Range("E33") "  -2.45%  "

# Row 34
Set-TextCell $ws.Range("E34") "  -4.84%  "

# Row 35
Set-TextCell $ws.Range("E35") "  -3.63%  "

# Row 36
Set-TextCell $ws.Range("E36") "  -0.12%  "

# Row 37
Set-TextCell $ws.Range("D37") "1.56"
Set-TextCell $ws.Range("E37") "  -5.00%  "

# Row 38
Set-TextCell $ws.Range("D38") "160.86"
Set-TextCell $ws.Range("E38") "  -1.63%  "

# Row 39
Set-TextCell $ws.Range("E39") "  +0.75%  "

# Row 40
Set-TextCell $ws.Range("D40") "28.18"
Set-TextCell $ws.Range("E40") "  +5.89%  "

# Row 41
Set-TextCell $ws.Range("E41") "  -4.46%  "

# Row 42
Set-TextCell $ws.Range("E42") "  -4.02%  "

# Row 43
Set-TextCell $ws.Range("D43") "6.65"
Set-TextCell $ws.Range("E43") "  -5.17%  "

# Row 44
Set-TextCell $ws.Range("E44") "  -4.26%  "

# Row 45
Set-TextCell $ws.Range("E45") "  -3.73%  "

# Row 46
Set-TextCell $ws.Range("D46") "2.720.04"
Set-TextCell $ws.Range("E46") "  -6.45%  "

# Row 47
Set-TextCell $ws.Range("D47") "25.94"
Set-TextCell $ws.Range("E47") "  -4.03%  "

# Row 48
Set-TextCell $ws.Range("E48") "  -2.65%  "

# Row 49
Set-TextCell $ws.Range("E49") "  -2.89%  "

# Row 50
Set-TextCell $ws.Range("D50") "329.36"
Set-TextCell $ws.Range("E50") "  -6.41%  "

# Row 51
Set-TextCell $ws.Range("D51") "1.05"
Set-TextCell $ws.Range("E51") "  -3.42%  "
